# workers_rank_mat.xlsx — refresh the ranking table with a re-run of the
# underlying simulation: per-worker "matrices" scores got new random draws,
# which re-shuffles a few same-gender ties (index/prolificid/name/gender)
# within each race group and updates every "matrices" score.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Asian group (rows 2-13): row 3/4/5 identities rotate, scores refresh ---
$ws.Range("G2").Value = 14.36475064273752

$ws.Range("C3").Value = 2
$ws.Range("D3").Value = "5f2c1a97a6809c060fec8820"
$ws.Range("E3").Value = "Maggie"
$ws.Range("G3").Value = 13.4427811560038

$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "60b1742bce2b39e0f1d19a1a"
$ws.Range("E4").Value = "Sabrina"
$ws.Range("G4").Value = 13.32257368402617

$ws.Range("C5").Value = 3
$ws.Range("D5").Value = "60bd88b8fc436774352f53b9"
$ws.Range("E5").Value = "Annes"
$ws.Range("G5").Value = 13.02548504840682

$ws.Range("G6").Value = 12.16366162123603
$ws.Range("G7").Value = 10.35758251781631
$ws.Range("G8").Value = 10.23661900101856
$ws.Range("G9").Value = 8.201924197465678
$ws.Range("G10").Value = 5.441561929436489
$ws.Range("G11").Value = 2.330660576781288
$ws.Range("G12").Value = 2.005372734962068
$ws.Range("G13").Value = 1.34066941120993

# --- Hispanic group (rows 14-25): rows 16/17 and 20/21/22 swap identities ---
$ws.Range("G14").Value = 15.02328293437414
$ws.Range("G15").Value = 11.12005548300506

$ws.Range("C16").Value = 3
$ws.Range("D16").Value = "60ba8ba51a5e0a105396888a"
$ws.Range("E16").Value = "Alfredo"
$ws.Range("F16").Value = "male"
$ws.Range("G16").Value = 10.3560449567461

$ws.Range("C17").Value = 2
$ws.Range("D17").Value = "60778ed0fde3e9c3a96f1d11"
$ws.Range("E17").Value = "Melissa"
$ws.Range("F17").Value = "female"
$ws.Range("G17").Value = 10.15590669353794

$ws.Range("G18").Value = 9.075645813370125
$ws.Range("G19").Value = 7.229575176107406

$ws.Range("C20").Value = 7
$ws.Range("D20").Value = "6024c18b094ac71dd93f4f5a"
$ws.Range("E20").Value = "Katherine"
$ws.Range("G20").Value = 5.144726965691964

$ws.Range("C21").Value = 9
$ws.Range("D21").Value = "5e35d91ea42bce592e996843"
$ws.Range("E21").Value = "Sergio"
$ws.Range("F21").Value = "male"
$ws.Range("G21").Value = 5.106254872490608

$ws.Range("C22").Value = 8
$ws.Range("D22").Value = "5f0142aa1eb1e528e7abce50"
$ws.Range("E22").Value = "Valeria"
$ws.Range("F22").Value = "female"
$ws.Range("G22").Value = 5.051234491524045

$ws.Range("G23").Value = 4.078136080597864
$ws.Range("G24").Value = 3.427904729701768
$ws.Range("G25").Value = 3.301880844181574
